$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New section header (row 16) ---
$ws.Range("B16").Value = "Project Number"

# --- Second capital-budgeting table header row (row 18) ---
$ws.Range("A18").Value = "Year"
$ws.Range("B18").Value = "Initial Investment"

$cols = @("D","E","F","G","H","I","J","K")
foreach ($col in $cols) {
  $ws.Range($col + "18").Value = -2000
}

# --- Data rows 19-33: project numbers in column B, cash flows in D:K ---
$data = @{
  19 = @{ B = 1;  D = 330;  E = 1666; F = 0; G = 160;  H = 280; I = 2200; J = 1200; K = -350 }
  20 = @{ B = 2;  D = 330;  E = 334;  F = 0; G = 200;  H = 280;           J = 900;  K = -60 }
  21 = @{ B = 3;  D = 330;  E = 165;  F = 0; G = 350;  H = 280;           J = 300;  K = 60 }
  22 = @{ B = 4;  D = 330;            F = 0; G = 395;  H = 280;           J = 90;   K = 350 }
  23 = @{ B = 5;  D = 330;            F = 0; G = 432;  H = 280;           J = 70;   K = 700 }
  24 = @{ B = 6;  D = 330;            F = 0; G = 440;  H = 280;                     K = 1200 }
  25 = @{ B = 7;  D = 330;            F = 0; G = 442;  H = 280;                     K = 2250 }
  26 = @{ B = 8;  D = 1000;           F = 0; G = 444;  H = 280 }
  27 = @{ B = 9;                      F = 0; G = 446;  H = 280 }
  28 = @{ B = 10;                     F = 0; G = 448;  H = 280 }
  29 = @{ B = 11;                     F = 0; G = 450;  H = 280 }
  30 = @{ B = 12;                     F = 0; G = 451;  H = 280 }
  31 = @{ B = 13;                     F = 0; G = 451;  H = 280 }
  32 = @{ B = 14;                     F = 0; G = 452;  H = 280 }
  33 = @{ B = 15;                     F = 10000; G = -2000; H = 280 }
}

foreach ($row in $data.Keys) {
  $rowData = $data[$row]
  foreach ($col in $rowData.Keys) {
    $ws.Range("$col$row").Value = $rowData[$col]
  }
}

# --- WACC row (35) ---
$ws.Range("B35").Value = "WACC"
$ws.Range("C35").Value = 0.1
$ws.Range("C35").NumberFormat = "0%"

# --- IRR row (36), D36:K36 as a shared formula (relative column references) ---
$ws.Range("B36").Value = "IRR"
$ws.Range("D36:K36").Formula = "=IRR(D18:D33)"
$ws.Range("D36:K36").NumberFormat = "0.00%"

# --- NPV row (37), D37:K37 as a shared formula (relative column references) ---
$ws.Range("B37").Value = "NPV"
$ws.Range("D37:K37").Style = "Currency"
$ws.Range("D37:K37").Formula = "=NPV(`$C`$35,D19:D33)+D18"

# --- Column D is now wider to fit "Initial Investment" spilling from column B's header context ---
$ws.Range("D1").EntireColumn.ColumnWidth = 20.42578125
